# Slide 10 ("3 aggregation methods tested" / model results slide): the
# "TextBox 10" shape was dragged so it now renders on top of (and is
# serialized after) the title, content placeholder and the four result
# pictures instead of being the first shape in the slide's shape tree.
#
# Bringing the shape to the front moves it to the end of the z-order
# (and therefore to the end of p:spTree), which is exactly the
# reordering the diff shows - no geometry/text inside the shape itself
# changes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

$textBox = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "TextBox 10") {
        $textBox = $shape
        break
    }
}

if ($textBox -ne $null) {
    $textBox.ZOrder(0)  # msoBringToFront -> moves shape to the end of spTree
}
